$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Liverpool", 12, 5, 4, 0, 1, 10, 1, 9),
    @("City", 12, 4, 4, 0, 0, 11, 3, 8),
    @("Aston Villa", 12, 5, 4, 0, 1, 10, 7, 3),
    @("Chelsea", 10, 5, 3, 1, 1, 11, 5, 6),
    @("Arsenal", 10, 4, 3, 1, 0, 6, 1, 5),
    @("Newcastle", 10, 5, 3, 1, 1, 7, 6, 1),
    @("Brighton", 8, 4, 2, 2, 0, 6, 2, 4),
    @("Fulham", 8, 5, 2, 2, 1, 7, 5, 2),
    @("Nottingham Forest", 8, 4, 2, 2, 0, 4, 2, 2),
    @("Tottenham", 7, 5, 2, 1, 2, 9, 5, 4),
    @("Manchester United", 7, 5, 2, 1, 2, 5, 5, 0),
    @("Brentford", 6, 5, 2, 0, 3, 7, 9, -2),
    @("Bournemouth", 5, 5, 1, 2, 2, 5, 8, -3),
    @("West Ham", 4, 5, 1, 1, 3, 5, 9, -4),
    @("Leicester City", 3, 5, 0, 3, 2, 6, 8, -2),
    @("Crystal Palace", 3, 5, 0, 3, 2, 4, 7, -3),
    @("Ipswich Town", 3, 5, 0, 3, 2, 3, 8, -5),
    @("Southampton", 1, 5, 0, 1, 4, 2, 9, -7),
    @("Everton", 1, 5, 0, 1, 4, 5, 14, -9),
    @("Wolves", 1, 5, 0, 1, 4, 5, 14, -9)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $colNum = $j + 1
        $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$j]
    }
}
